# Add a new data row (row 2) to track 12:30-to-close return data.
# Dates and percent-like strings must be forced to text so they are
# stored verbatim (matching "06/29/23" and "0.10282%") instead of being
# auto-converted to a date serial / percentage fraction by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: Trading Date (text, not a real date)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "06/29/23"
$ws.Range("A2").NumberFormat = "General"
$ws.Range("A2").Style = "Normal"

# B2: Gamma Imbalance (plain number)
$ws.Range("B2").Value = 0.0256
$ws.Range("B2").Style = "Normal"

# C2: Hedge Pressure (plain number)
$ws.Range("C2").Value = 0.00763
$ws.Range("C2").Style = "Normal"

# D2: 12:30 to close return (text, keeps the literal "%" as typed)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.10282%"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"

# E2 (t+1 first 30 return) intentionally left blank - not yet known.
